$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# "dbExcel"/"WebExcel" columns (and their data) one column to the right.
$ws.Columns("B").Insert()

# Set the new column's width to match column A (75.81640625 characters,
# stored in character units). The runtime quantizes ColumnWidth to the
# nearest 1/6 character, so 75 is the closest input that reproduces
# column A's effective stored width.
$ws.Columns("B").ColumnWidth = 75

# New header for the inserted column
$ws.Range("B1").Value = "StatQuery"

# New query text for the inserted column, matching the style (wrap text) of A2
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lip and oropharyngeal neoplasms malignant :: Melanoma-mucosa/maxilla']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").WrapText = $true

# Update the selection to just B2 (as in the edited workbook)
$ws.Range("B2").Select()
